# "Generate Report for Archive"
#
# The localization run moved on: the status that used to read
# "Ready for handoff" is now "In Translation" everywhere it appears
# (Overview!E2:F2, zh-cn!C2, de-de!C2 - all the same shared string).
# Because the new status text is shorter, the Status-ish columns that
# host it are re-sized down to fit.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# NOTE on the column width below: this host snaps COM `ColumnWidth`
# writes to Excel's internal 1/6-character pixel grid
# (stored = ROUND(ColumnWidth*6 + 5) / 6), same as real Excel COM
# automation. Feeding it 12.5 lands on the nearest attainable grid
# point to the archived report's target width, same family as the
# narrower 13.41-char column the generator produced.
$narrowColumnWidth = 12.5

# --- Overview sheet: zh-cn (E) and de-de (F) status cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $narrowColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $narrowColumnWidth

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $narrowColumnWidth

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $narrowColumnWidth
